$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Values
$ws.Range("B1").Value = 0
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "disconnected_elements"

# Style B1: bold font, thin box border, centered/top-aligned
$r1 = $ws.Range("B1")
$r1.Font.Bold = $true
$r1.Borders.LineStyle = 1
$r1.Borders.Weight = 2
$r1.VerticalAlignment = -4160
$r1.HorizontalAlignment = -4108

# Copy the same formatting onto A2 so both cells share one cellXf
$r1.Copy()
$r2 = $ws.Range("A2")
$r2.PasteSpecial(-4122)
